$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 177 ("Hortaliza, Vega
# Monumental Concepción - Acelga", weekly refresh). Excel shifts the
# existing rows 177:296 down to 178:297 and grows the used range to
# A1:R297, matching the canonical diff.
$ws.Rows("177:177").Insert()

# Populate the newly inserted row 177 with the new weekly record.
$ws.Range("A177").Value = 11
$ws.Range("B177").Value = "Vega Monumental Concepción"
$ws.Range("C177").Value = "Bíobío"
$ws.Range("D177").Value = 44806
$ws.Range("E177").Value = 8
$ws.Range("F177").Value = 100112009
$ws.Range("G177").Value = "Acelga"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 350
$ws.Range("K177").Value = 700
$ws.Range("L177").Value = 750
$ws.Range("M177").Value = 729
$ws.Range("N177").Value = "`$/atado"
$ws.Range("O177").Value = "Región de Ñuble"
$ws.Range("P177").Value = 729
$ws.Range("Q177").Value = 1
$ws.Range("R177").Value = "Hortaliza"
